# Weekly update: add the latest week's "Ajo" (garlic) price records for
# "Vega Central Mapocho de Santiago" by inserting two new rows at the top
# of this sub-group (row 173), pushing the existing rows 173-180 down to
# 175-182.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 173 (shifts 173:180 -> 175:182)
$ws.Rows("173:174").Insert()

# New row 173: Ajo / Chino / Primera
$ws.Cells.Item(173, 1).Value = 9
$ws.Cells.Item(173, 2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(173, 3).Value = 'Metropolitana'
$ws.Cells.Item(173, 4).Value = 44585
$ws.Cells.Item(173, 5).Value = 13
$ws.Cells.Item(173, 6).Value = 100112003
$ws.Cells.Item(173, 7).Value = 'Ajo'
$ws.Cells.Item(173, 8).Value = 'Chino'
$ws.Cells.Item(173, 9).Value = 'Primera'
$ws.Cells.Item(173, 10).Value = 610
$ws.Cells.Item(173, 11).Value = 17500
$ws.Cells.Item(173, 12).Value = 18000
$ws.Cells.Item(173, 13).Value = 17750
$ws.Cells.Item(173, 14).Value = '$/caja 10 kilos'
$ws.Cells.Item(173, 15).Value = 'China'
$ws.Cells.Item(173, 16).Value = 1775
$ws.Cells.Item(173, 17).Value = 10
$ws.Cells.Item(173, 18).Value = 'Hortaliza'

# New row 174: Ajo / Chino / Primera
$ws.Cells.Item(174, 1).Value = 9
$ws.Cells.Item(174, 2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(174, 3).Value = 'Metropolitana'
$ws.Cells.Item(174, 4).Value = 44585
$ws.Cells.Item(174, 5).Value = 13
$ws.Cells.Item(174, 6).Value = 100112003
$ws.Cells.Item(174, 7).Value = 'Ajo'
$ws.Cells.Item(174, 8).Value = 'Chino'
$ws.Cells.Item(174, 9).Value = 'Primera'
$ws.Cells.Item(174, 10).Value = 340
$ws.Cells.Item(174, 11).Value = 18000
$ws.Cells.Item(174, 12).Value = 18500
$ws.Cells.Item(174, 13).Value = 18250
$ws.Cells.Item(174, 14).Value = '$/malla 10 kilos'
$ws.Cells.Item(174, 15).Value = 'China'
$ws.Cells.Item(174, 16).Value = 1825
$ws.Cells.Item(174, 17).Value = 10
$ws.Cells.Item(174, 18).Value = 'Hortaliza'
